# Update "想去人数" (want-to-go count, column F) figures scraped for this
# gh-pages data refresh. Two sheets carry the same event list:
#   "展览"   (Worksheets index 1) - rows 2-32
#   "全部类型" (Worksheets index 4) - rows 2-33 (one extra event, shifted by 1 row)

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 191
    3  = 5378
    8  = 583
    9  = 1053
    11 = 1474
    12 = 4348
    13 = 441
    17 = 3458
    18 = 170
    19 = 1097
    23 = 127
    26 = 73
    27 = 313
    32 = 28
}
foreach ($row in $sheet1Updates.Keys) {
    $sheet1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 191
    4  = 5378
    9  = 583
    10 = 1053
    12 = 1474
    13 = 4348
    14 = 441
    18 = 3458
    19 = 170
    20 = 1097
    24 = 127
    27 = 73
    28 = 313
    33 = 28
}
foreach ($row in $sheet4Updates.Keys) {
    $sheet4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
